$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") from 2023-10-03 (45202) to 2023-10-04 (45203)
# for all data rows (row 2 through row 163).
$lastRow = 163
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45203
}
